$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking Price cells as Text so the literal digits survive
# (values such as "1.000" / "316.06" would otherwise be parsed as numbers).
$forceTextCells = @("D4","D5","D6","D7","D8","D11","D12","D13","D14","D16","D18","D19","D20","D21","D22","D23","D25","D26","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Per-row cell updates -------------------------------------------------
# Row 2
$ws.Range("D2").Value = "24.514.47"
$ws.Range("E2").Value = "  -0.98%  "

# Row 3
$ws.Range("D3").Value = "1.695.03"
$ws.Range("E3").Value = "  -0.47%  "

# Row 4
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.32%  "

# Row 5
$ws.Range("D5").Value = "316.06"
$ws.Range("E5").Value = "  -0.32%  "

# Row 6
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.38%  "

# Row 7
$ws.Range("D7").Value = "0.3922"
$ws.Range("E7").Value = "  -0.70%  "

# Row 8
$ws.Range("D8").Value = "0.4069"
$ws.Range("E8").Value = "  +0.42%  "

# Row 9
$ws.Range("E9").Value = "  -2.73%  "

# Row 10
$ws.Range("E10").Value = "  -0.31%  "

# Row 11
$ws.Range("D11").Value = "52.49"
$ws.Range("E11").Value = "  -2.17%  "

# Row 12
$ws.Range("D12").Value = "0.08789"
$ws.Range("E12").Value = "  -1.19%  "

# Row 13
$ws.Range("D13").Value = "26.63"
$ws.Range("E13").Value = "  +12.09%  "

# Row 14
$ws.Range("D14").Value = "7.511"
$ws.Range("E14").Value = "  +0.39%  "

# Row 15
$ws.Range("E15").Value = "  -1.73%  "

# Row 16
$ws.Range("D16").Value = "0.00001350"
$ws.Range("E16").Value = "  +1.69%  "

# Row 17
$ws.Range("D17").Value = "1.686.73"
$ws.Range("E17").Value = "  -1.06%  "

# Row 18
$ws.Range("D18").Value = "98.08"
$ws.Range("E18").Value = "  -1.82%  "

# Row 19
$ws.Range("D19").Value = "0.07151"
$ws.Range("E19").Value = "  +1.23%  "

# Row 20
$ws.Range("D20").Value = "20.58"
$ws.Range("E20").Value = "  +4.02%  "

# Row 21
$ws.Range("D21").Value = "7.287"
$ws.Range("E21").Value = "  +2.78%  "

# Row 22
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.20%  "

# Row 23
$ws.Range("D23").Value = "14.34"
$ws.Range("E23").Value = "  -2.12%  "

# Row 24
$ws.Range("D24").Value = "24.503.27"
$ws.Range("E24").Value = "  -0.96%  "

# Row 25
$ws.Range("D25").Value = "3.010"
$ws.Range("E25").Value = "  -6.49%  "

# Row 26
$ws.Range("D26").Value = "2.323"
$ws.Range("E26").Value = "  -2.15%  "

# Row 27
$ws.Range("E27").Value = "  -0.47%  "

# Row 28
$ws.Range("D28").Value = "166.76"
$ws.Range("E28").Value = "  +2.45%  "

# Row 29
$ws.Range("D29").Value = "8.543"
$ws.Range("E29").Value = "  -2.86%  "

# Row 30
$ws.Range("D30").Value = "5.397"
$ws.Range("E30").Value = "  +4.29%  "

# Row 31
$ws.Range("D31").Value = "139.35"
$ws.Range("E31").Value = "  +2.40%  "

# Row 32
$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "1.871.53"
$ws.Range("E32").Value = "  -1.26%  "

# Row 33
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "2.195"
$ws.Range("E33").Value = "  +10.12%  "

# Row 34
$ws.Range("D34").Value = "0.08774"
$ws.Range("E34").Value = "  -2.34%  "

# Row 35
$ws.Range("D35").Value = "7.330"
$ws.Range("E35").Value = "  -5.08%  "

# Row 36
$ws.Range("D36").Value = "1.038"
$ws.Range("E36").Value = "  -4.05%  "

# Row 37
$ws.Range("D37").Value = "0.02992"
$ws.Range("E37").Value = "  +7.15%  "

# Row 38
$ws.Range("D38").Value = "0.2784"
$ws.Range("E38").Value = "  +0.71%  "

# Row 39
$ws.Range("D39").Value = "10.95"
$ws.Range("E39").Value = "  -1.59%  "

# Row 40
$ws.Range("D40").Value = "0.09159"
$ws.Range("E40").Value = "  -0.30%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "14.20"
$ws.Range("E41").Value = "  -2.65%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.8041"
$ws.Range("E42").Value = "  +4.04%  "

# Row 43
$ws.Range("D43").Value = "1.477"
$ws.Range("E43").Value = "  +1.03%  "

# Row 44
$ws.Range("D44").Value = "17.58"
$ws.Range("E44").Value = "  +10.07%  "

# Row 45
$ws.Range("D45").Value = "2.667"
$ws.Range("E45").Value = "  +3.27%  "

# Row 46
$ws.Range("D46").Value = "0.7275"
$ws.Range("E46").Value = "  +0.73%  "

# Row 47
$ws.Range("D47").Value = "4.259"
$ws.Range("E47").Value = "  +0.83%  "

# Row 48
$ws.Range("D48").Value = "1.405"
$ws.Range("E48").Value = "  +3.23%  "

# Row 49
$ws.Range("D49").Value = "0.9998"
$ws.Range("E49").Value = "  -0.59%  "

# Row 50
$ws.Range("D50").Value = "141.07"
$ws.Range("E50").Value = "  +0.18%  "

# Row 51
$ws.Range("D51").Value = "0.08165"
$ws.Range("E51").Value = "  +2.17%  "

# Restore default (Normal) style on the forced-text cells so no stray number
# format is left attached to the cell (matches original unstyled cells).
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).ClearFormats()
}

